$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 191 - 04-10-2021
$a191 = $ws.Cells.Item(191, 1)
$a191.NumberFormat = "@"
$a191.Value = "04-10-2021"
$a191.Style = "Normal"
$ws.Cells.Item(191, 2).Value = 0.15
$ws.Cells.Item(191, 3).Value = 0.21
$ws.Cells.Item(191, 4).Value = 0.01
$ws.Cells.Item(191, 5).Value = 0.35
$ws.Cells.Item(191, 6).Value = 0.5

# Row 192 - 05-10-2021
$a192 = $ws.Cells.Item(192, 1)
$a192.NumberFormat = "@"
$a192.Value = "05-10-2021"
$a192.Style = "Normal"
$ws.Cells.Item(192, 2).Value = 0.16
$ws.Cells.Item(192, 3).Value = 0.25
$ws.Cells.Item(192, 4).Value = 0.01
$ws.Cells.Item(192, 5).Value = 0.15
$ws.Cells.Item(192, 6).Value = 0.47

# Row 193 - 06-10-2021
$a193 = $ws.Cells.Item(193, 1)
$a193.NumberFormat = "@"
$a193.Value = "06-10-2021"
$a193.Style = "Normal"
$ws.Cells.Item(193, 2).Value = 0.16
$ws.Cells.Item(193, 3).Value = 0.23
$ws.Cells.Item(193, 4).Value = 0
$ws.Cells.Item(193, 5).Value = 0.17
$ws.Cells.Item(193, 6).Value = 0.35
